$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price observation (row) is inserted in the middle of the table at
# row 323. Excel's native row-insert semantics push every existing row at
# and below 323 down by one (so the former row 323 becomes 324, ..., the
# former row 396 becomes 397), carrying each row's formatting along with
# it (matching the surrounding date/number styling).
$ws.Rows("323:323").Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A323").Value2 = 4
$ws.Range("B323").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C323").Value2 = "Los Lagos"
$ws.Range("D323").Value2 = 45015
$ws.Range("E323").Value2 = 10
$ws.Range("F323").Value2 = 100112021
$ws.Range("G323").Value2 = "Ají"
$ws.Range("H323").Value2 = "Inferno"
$ws.Range("I323").Value2 = "Primera"
$ws.Range("J323").Value2 = 50
$ws.Range("K323").Value2 = 28000
$ws.Range("L323").Value2 = 28000
$ws.Range("M323").Value2 = 28000
$ws.Range("N323").Value2 = "`$/caja 15 kilos"
$ws.Range("O323").Value2 = "Provincia de Quillota"
$ws.Range("P323").Value2 = 1867
$ws.Range("Q323").Value2 = 15
$ws.Range("R323").Value2 = "Hortaliza"
